$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 24, shifting existing rows 24-45 down to 26-47.
$ws.Rows("24:25").Insert()

# New row 24: Tuna, Especial quality, 44601 (2022-02-09), $/caja 18 kilos
$ws.Range("A24").Value2 = 10
$ws.Range("B24").Value2 = "Vega Modelo de Temuco"
$ws.Range("C24").Value2 = "La Araucanía"
$ws.Range("D24").Value2 = 44601
$ws.Range("E24").Value2 = 9
$ws.Range("F24").Value2 = "Fruta"
$ws.Range("G24").Value2 = 100107
$ws.Range("H24").Value2 = "Otros"
$ws.Range("I24").Value2 = 100107011
$ws.Range("J24").Value2 = "Tuna"
$ws.Range("K24").Value2 = "Sin especificar"
$ws.Range("L24").Value2 = "Especial"
$ws.Range("M24").Value2 = 30
$ws.Range("N24").Value2 = 25000
$ws.Range("O24").Value2 = 25000
$ws.Range("P24").Value2 = 25000
$ws.Range("Q24").Value2 = "$/caja 18 kilos"
$ws.Range("R24").Value2 = "Provincia de Los Andes"
$ws.Range("S24").Value2 = 1389
$ws.Range("T24").Value2 = 18

# New row 25: Tuna, Primera quality, 44601 (2022-02-09), $/caja 18 kilos
$ws.Range("A25").Value2 = 10
$ws.Range("B25").Value2 = "Vega Modelo de Temuco"
$ws.Range("C25").Value2 = "La Araucanía"
$ws.Range("D25").Value2 = 44601
$ws.Range("E25").Value2 = 9
$ws.Range("F25").Value2 = "Fruta"
$ws.Range("G25").Value2 = 100107
$ws.Range("H25").Value2 = "Otros"
$ws.Range("I25").Value2 = 100107011
$ws.Range("J25").Value2 = "Tuna"
$ws.Range("K25").Value2 = "Sin especificar"
$ws.Range("L25").Value2 = "Primera"
$ws.Range("M25").Value2 = 80
$ws.Range("N25").Value2 = 18000
$ws.Range("O25").Value2 = 18000
$ws.Range("P25").Value2 = 18000
$ws.Range("Q25").Value2 = "$/caja 18 kilos"
$ws.Range("R25").Value2 = "Provincia de Los Andes"
$ws.Range("S25").Value2 = 1000
$ws.Range("T25").Value2 = 18
